# Update crypto price (D) and volume-change (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.980.43"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "2.565.90"
$ws.Range("E3").Value = "  -3.08%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'516.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'142.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "2.580.82"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").Value = "'6.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("E12").Value = "  -4.66%  "
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").Value = "3.019.80"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "57.991.46"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "'20.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "2.578.25"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").Value = "'342.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "'10.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").Value = "'6.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'65.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  -5.93%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "2.690.56"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("D29").Value = "'6.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("D30").Value = "0.0₃0743"
$ws.Range("E30").Value = "  -7.12%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'6.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.04%  "
$ws.Range("D33").Value = "'1.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D35").Value = "'149.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").Value = "'3.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.02%  "
$ws.Range("E37").Value = "  -4.06%  "
$ws.Range("D38").Value = "'0.870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("D39").Value = "'36.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").Value = "'1.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "'0.831"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "'269.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "'10.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D47").Value = "'0.585"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("D48").Value = "'18.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").Value = "1.976.82"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "'4.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.77%  "
